# Swap the contents of specific columns between row 7 and row 8.
# Columns affected: A, I, P, Q, R, Z, AB, AC

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "I", "P", "Q", "R", "Z", "AB", "AC")

foreach ($col in $cols) {
    $cell7 = $ws.Range($col + "7")
    $cell8 = $ws.Range($col + "8")

    $val7 = $cell7.Value2
    $val8 = $cell8.Value2

    $cell7.Value2 = $val8
    $cell8.Value2 = $val7
}
